# Completes the testing matrix: rewrites column A with the updated set of
# test-case questions (some rows reordered/renamed, many new rows appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Functioneaza citirea din fisiere?",
    "Functioneaza calcularea error screen?",
    "Functioneaza calcularea de screen skips?",
    "Se poate da skip unui pas?",
    "Daca apasam enter, trecem la urmatorul pas?",
    "Functia userInteraction functioneaza?",
    "Functioneaza functia deleteFlow?",
    "Functioneaza functia runFlow?",
    "Putem oferi un nume flow-ului?",
    "Putem oferi de la tastatura detalii legate de flow?",
    "Putem alege tipul de step pe care dorim sa-l executam?",
    "Putem alege tipul de operatie pe care vrem sa-l executam in Calculus Step?",
    "Pot alege orice alta litera in afara de n pentru a continua?",
    "Functia getTypeStep returneaza numele step-ului pe care eu il vrem?",
    "La final flow-ul este considerat complet?",
    "Putem oferi valori pentru fiecare parametru al fiecarui step?",
    "Daca tipul de step ales este Output, primim informatii de la celelalte step-uri executate?",
    "La final primim informatii despre ce step-uri am completat in timpul rularii flow-ului?",
    "Atunci cand rularea flow-ului se incheie, putem afla momentul in care a fost creat flow-ul?",
    "Putem oferi o descriere pentru un flow?"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A20").Select()
